{"js": "// Apply the text replacements described by the diff: the date line and the\n// 25 two-digit-multiplication answer cells in the table.\nconst replacements = [\n  [\"2023-11-22 Wednesday\", \"2023-11-23 Thursday\"],\n  [\"54\u00d737=1998\", \"11\u00d752=572\"],\n  [\"27\u00d740=1080\", \"57\u00d798=5586\"],\n  [\"48\u00d719=912\", \"84\u00d754=4536\"],\n  [\"34\u00d741=1394\", \"24\u00d762=1488\"],\n  [\"39\u00d780=3120\", \"49\u00d722=1078\"],\n  [\"17\u00d788=1496\", \"50\u00d724=1200\"],\n  [\"35\u00d796=3360\", \"43\u00d714=602\"],\n  [\"43\u00d779=3397\", \"21\u00d736=756\"],\n  [\"93\u00d719=1767\", \"94\u00d797=9118\"],\n  [\"83\u00d745=3735\", \"42\u00d747=1974\"],\n  [\"98\u00d793=9114\", \"34\u00d781=2754\"],\n  [\"23\u00d767=1541\", \"26\u00d721=546\"],\n  [\"31\u00d779=2449\", \"44\u00d725=1100\"],\n  [\"15\u00d774=1110\", \"61\u00d757=3477\"],\n  [\"65\u00d740=2600\", \"29\u00d724=696\"],\n  [\"75\u00d722=1650\", \"55\u00d756=3080\"],\n  [\"34\u00d760=2040\", \"71\u00d735=2485\"],\n  [\"31\u00d751=1581\", \"47\u00d791=4277\"],\n  [\"37\u00d731=1147\", \"71\u00d782=5822\"],\n  [\"30\u00d712=360\", \"65\u00d732=2080\"],\n  [\"12\u00d712=144\", \"75\u00d712=900\"],\n  [\"99\u00d794=9306\", \"55\u00d723=1265\"],\n  [\"46\u00d768=3128\", \"59\u00d764=3776\"],\n  [\"98\u00d774=7252\", \"95\u00d749=4655\"],\n  [\"94\u00d766=6204\", \"47\u00d764=3008\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Apply the text replacements described by the diff: the date line and the\n# 25 two-digit-multiplication answer cells in the table.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2023-11-22 Wednesday\", \"2023-11-23 Thursday\"),\n    @(\"54\u00d737=1998\", \"11\u00d752=572\"),\n    @(\"27\u00d740=1080\", \"57\u00d798=5586\"),\n    @(\"48\u00d719=912\", \"84\u00d754=4536\"),\n    @(\"34\u00d741=1394\", \"24\u00d762=1488\"),\n    @(\"39\u00d780=3120\", \"49\u00d722=1078\"),\n    @(\"17\u00d788=1496\", \"50\u00d724=1200\"),\n    @(\"35\u00d796=3360\", \"43\u00d714=602\"),\n    @(\"43\u00d779=3397\", \"21\u00d736=756\"),\n    @(\"93\u00d719=1767\", \"94\u00d797=9118\"),\n    @(\"83\u00d745=3735\", \"42\u00d747=1974\"),\n    @(\"98\u00d793=9114\", \"34\u00d781=2754\"),\n    @(\"23\u00d767=1541\", \"26\u00d721=546\"),\n    @(\"31\u00d779=2449\", \"44\u00d725=1100\"),\n    @(\"15\u00d774=1110\", \"61\u00d757=3477\"),\n    @(\"65\u00d740=2600\", \"29\u00d724=696\"),\n    @(\"75\u00d722=1650\", \"55\u00d756=3080\"),\n    @(\"34\u00d760=2040\", \"71\u00d735=2485\"),\n    @(\"31\u00d751=1581\", \"47\u00d791=4277\"),\n    @(\"37\u00d731=1147\", \"71\u00d782=5822\"),\n    @(\"30\u00d712=360\", \"65\u00d732=2080\"),\n    @(\"12\u00d712=144\", \"75\u00d712=900\"),\n    @(\"99\u00d794=9306\", \"55\u00d723=1265\"),\n    @(\"46\u00d768=3128\", \"59\u00d764=3776\"),\n    @(\"98\u00d774=7252\", \"95\u00d749=4655\"),\n    @(\"94\u00d766=6204\", \"47\u00d764=3008\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $d.Content.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
